# Natmi following Dr Hou advice
# Rebuild the LR-pair data block (rows 2-10) on Sheet1 with the updated
# sender/receiver cluster matrix (adds the "sCs" cluster) and refreshed values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row number -> full A:T row values, in the final target layout
$data = @{
  2 = @("ECs", "Il23a", "Il12rb1", "ECs", 3, 1, 10.81092366666667, 32.432771, 0.886286770790134, 0.8862867707901342, 1, 0.3333333333333333, 0.157372, 0.472116, 0.1623752133230521, 0.1623752133230521, 1.701336679270667, 15.312030113436, 0.143911003472447, 0.143911003472447)
  3 = @("ECs", "Il23a", "Il12rb1", "FAPs", 3, 1, 10.81092366666667, 32.432771, 0.886286770790134, 0.8862867707901342, 3, 1, 0.6279076666666666, 1.883723, 0.6478702775727568, 0.6478702775727568, 6.788261854048112, 61.094356686433, 0.5741988562008664, 0.5741988562008665)
  4 = @("ECs", "Il23a", "Il12rb1", "sCs", 3, 1, 10.81092366666667, 32.432771, 0.886286770790134, 0.8862867707901342, 3, 1, 0.1839076666666667, 0.5517230000000001, 0.1897545091041911, 0.1897545091041911, 1.988211746048111, 17.893905714433, 0.1681769111168206, 0.1681769111168206)
  5 = @("FAPs", "Il23a", "Il12rb1", "ECs", 3, 1, 0.9380866666666666, 2.81426, 0.07690497390937834, 0.07690497390937835, 1, 0.3333333333333333, 0.157372, 0.472116, 0.1623752133230521, 0.1623752133230521, 0.1476285749066666, 1.32865717416, 0.01248746154413906, 0.01248746154413907)
  6 = @("FAPs", "Il23a", "Il12rb1", "FAPs", 3, 1, 0.9380866666666666, 2.81426, 0.07690497390937834, 0.07690497390937835, 3, 1, 0.6279076666666666, 1.883723, 0.6478702775727568, 0.6478702775727568, 0.5890318099977777, 5.30128628998, 0.04982444679339457, 0.04982444679339457)
  7 = @("FAPs", "Il23a", "Il12rb1", "sCs", 3, 1, 0.9380866666666666, 2.81426, 0.07690497390937834, 0.07690497390937835, 3, 1, 0.1839076666666667, 0.5517230000000001, 0.1897545091041911, 0.1897545091041911, 0.1725213299977778, 1.55269196998, 0.01459306557184471, 0.01459306557184471)
  8 = @("sCs", "Il23a", "Il12rb1", "ECs", 1, 0.3333333333333333, 0.448987, 1.346961, 0.03680825530048758, 0.03680825530048758, 1, 0.3333333333333333, 0.157372, 0.472116, 0.1623752133230521, 0.1623752133230521, 0.070657982164, 0.635921839476, 0.005976748306466033, 0.005976748306466034)
  9 = @("sCs", "Il23a", "Il12rb1", "FAPs", 1, 0.3333333333333333, 0.448987, 1.346961, 0.03680825530048758, 0.03680825530048758, 3, 1, 0.6279076666666666, 1.883723, 0.6478702775727568, 0.6478702775727568, 0.2819223795336667, 2.537301415803, 0.02384697457849579, 0.02384697457849579)
  10 = @("sCs", "Il23a", "Il12rb1", "sCs", 1, 0.3333333333333333, 0.448987, 1.346961, 0.03680825530048758, 0.03680825530048758, 3, 1, 0.1839076666666667, 0.5517230000000001, 0.1897545091041911, 0.1897545091041911, 0.08257215153366668, 0.7431493638030001, 0.006984532415525761, 0.00698453241552576)
}

foreach ($r in ($data.Keys | Sort-Object)) {
  $rowVals = $data[$r]
  for ($c = 0; $c -lt $rowVals.Length; $c++) {
    $ws.Cells.Item([int]$r, $c + 1).Value = $rowVals[$c]
  }
}
